# QA-7808 fixes for import cases and case search
#
# The "Cases" test-data fixture carried a stray caseid value (a leftover
# GUID) in A2 that no longer belongs in the sheet - the other columns in
# that row (name/owner info) are still valid. Clear it out and leave the
# selection sitting on A2, same as the fixture was left after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").ClearContents() | Out-Null
$ws.Range("A2").Select() | Out-Null
